# Apply the authored edits to PowerPointTest.pptx:
#  - Slide 1, "Table 1": Author cell value "Name" -> "Tara Keena"
#  - Slide 1, "Table 1": Version cell value "1.0" -> "1.1"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Tara Keena"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "1.1"
